$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.642.93'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.886.62'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.34'
$ws.Range('E5').Value = '  -4.44%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4878'
$ws.Range('E7').Value = '  -2.42%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2886'
$ws.Range('E8').Value = '  -4.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06658'
$ws.Range('E9').Value = '  -3.78%  '
$ws.Range('D10').Value = '1.880.65'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '16.77'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07234'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '88.76'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.000'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('D16').Value = '30.591.49'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000007836'
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.94'
$ws.Range('E19').Value = '  -3.90%  '
$ws.Range('D20').Value = '2.124.36'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.733'
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '186.97'
$ws.Range('E23').Value = '  +4.23%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.037'
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.262'
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.92'
$ws.Range('E26').Value = '  +3.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.23'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('E28').Value = '  -6.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.408'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.255'
$ws.Range('E30').Value = '  -2.74%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09027'
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.928'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05190'
$ws.Range('E33').Value = '  -1.50%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7313'
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('E35').Value = '  -6.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.696'
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01815'
$ws.Range('E37').Value = '  -5.63%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.656'
$ws.Range('E38').Value = '  -3.30%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9200'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.041'
$ws.Range('E40').Value = '  -7.45%  '
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '103.86'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9992'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.715'
$ws.Range('E44').Value = '  -4.06%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1339'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.270'
$ws.Range('E46').Value = '  -7.57%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05829'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.3953'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.617'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.407'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '33.19'
$ws.Range('E51').Value = '  -0.61%  '
